$wb = $excel.ActiveWorkbook

# Auto-generated edit script: updates market-price columns (H-N) for the
# rows touched by the scheduled price-refresh run, across 7 of the 8 leve sheets.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2268.8572
$ws.Range("I40").Value = 1150
$ws.Range("J40").Value = 2716.4
$ws.Range("K40").Value = 1150
$ws.Range("L40").Value = 2716.4
$ws.Range("M40").Value = -975
$ws.Range("N40").Value = -3066.4

$ws.Range("H53").Value = 287
$ws.Range("I53").Value = 226.36363
$ws.Range("J53").Value = 361.1111
$ws.Range("K53").Value = 226.36363
$ws.Range("L53").Value = 361.1111
$ws.Range("M53").Value = 410.63637
$ws.Range("N53").Value = -1635.1111

$ws.Range("H69").Value = 3903.4324
$ws.Range("I69").Value = 3918.3333
$ws.Range("J69").Value = 3900.5483
$ws.Range("K69").Value = 11754.9999
$ws.Range("L69").Value = 11701.6449
$ws.Range("M69").Value = -10880.9999
$ws.Range("N69").Value = -13449.6449

$ws.Range("H72").Value = 3903.4324
$ws.Range("I72").Value = 3918.3333
$ws.Range("J72").Value = 3900.5483
$ws.Range("K72").Value = 35264.9997
$ws.Range("L72").Value = 35104.9347
$ws.Range("M72").Value = -30896.9997
$ws.Range("N72").Value = -43840.9347

$ws.Range("H86").Value = 1816.238
$ws.Range("I86").Value = 1067.3334
$ws.Range("J86").Value = 2377.9167
$ws.Range("K86").Value = 1067.3334
$ws.Range("L86").Value = 2377.9167
$ws.Range("M86").Value = 55.66660000000002
$ws.Range("N86").Value = -4623.9167

$ws.Range("H89").Value = 1816.238
$ws.Range("I89").Value = 1067.3334
$ws.Range("J89").Value = 2377.9167
$ws.Range("K89").Value = 5336.666999999999
$ws.Range("L89").Value = 11889.5835
$ws.Range("M89").Value = 279.3330000000005
$ws.Range("N89").Value = -23121.5835

$ws.Range("H103").Value = 398.66666
$ws.Range("I103").Value = 385.5
$ws.Range("K103").Value = 1156.5
$ws.Range("M103").Value = -570.5

$ws.Range("H107").Value = 457.92307
$ws.Range("I107").Value = 462.75
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 462.75
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 1457.25
$ws.Range("N107").Value = -4240

$ws.Range("H135").Value = 415
$ws.Range("I135").Value = 327.66666
$ws.Range("J135").Value = 808
$ws.Range("K135").Value = 2948.99994
$ws.Range("L135").Value = 7272
$ws.Range("M135").Value = -413.9999399999997
$ws.Range("N135").Value = -12342

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3534.3125
$ws.Range("I86").Value = 2835.3333
$ws.Range("J86").Value = 3695.6155
$ws.Range("K86").Value = 2835.3333
$ws.Range("L86").Value = 3695.6155
$ws.Range("M86").Value = -1712.3333
$ws.Range("N86").Value = -5941.6155

$ws.Range("H89").Value = 3534.3125
$ws.Range("I89").Value = 2835.3333
$ws.Range("J89").Value = 3695.6155
$ws.Range("K89").Value = 14176.6665
$ws.Range("L89").Value = 18478.0775
$ws.Range("M89").Value = -8560.666499999999
$ws.Range("N89").Value = -29710.0775

$ws.Range("H107").Value = 2174.7827
$ws.Range("I107").Value = 1820.7142
$ws.Range("J107").Value = 2725.5557
$ws.Range("K107").Value = 1820.7142
$ws.Range("L107").Value = 2725.5557
$ws.Range("M107").Value = 99.28580000000011
$ws.Range("N107").Value = -6565.5557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1293.2307
$ws.Range("I122").Value = 1168
$ws.Range("J122").Value = 1575
$ws.Range("K122").Value = 3504
$ws.Range("L122").Value = 4725
$ws.Range("M122").Value = -1054
$ws.Range("N122").Value = -9625

$ws.Range("H132").Value = 27033860
$ws.Range("I132").Value = 50007748
$ws.Range("J132").Value = 5756.8823
$ws.Range("K132").Value = 150023244
$ws.Range("L132").Value = 17270.6469
$ws.Range("M132").Value = -150020714
$ws.Range("N132").Value = -22330.6469

$ws.Range("H134").Value = 2421.0588
$ws.Range("I134").Value = 2238.8572
$ws.Range("J134").Value = 3271.3333
$ws.Range("K134").Value = 6716.571599999999
$ws.Range("L134").Value = 9813.999899999999
$ws.Range("M134").Value = -4181.571599999999
$ws.Range("N134").Value = -14883.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 47.8
$ws.Range("I10").Value = 47.8
$ws.Range("K10").Value = 143.4
$ws.Range("M10").Value = -4.399999999999977

$ws.Range("H34").Value = 979.2273
$ws.Range("I34").Value = 272.75
$ws.Range("J34").Value = 1136.2222
$ws.Range("K34").Value = 818.25
$ws.Range("L34").Value = 3408.6666
$ws.Range("M34").Value = -734.25
$ws.Range("N34").Value = -3576.6666

$ws.Range("H39").Value = 3078.5715
$ws.Range("J39").Value = 3253.8462
$ws.Range("L39").Value = 9761.5386
$ws.Range("N39").Value = -10349.5386

$ws.Range("H55").Value = 2585
$ws.Range("I55").Value = 525
$ws.Range("J55").Value = 3100
$ws.Range("K55").Value = 1575
$ws.Range("L55").Value = 9300
$ws.Range("M55").Value = -1398
$ws.Range("N55").Value = -9654

$ws.Range("H130").Value = 2228.75
$ws.Range("I130").Value = 665
$ws.Range("J130").Value = 2750
$ws.Range("K130").Value = 1995
$ws.Range("L130").Value = 8250
$ws.Range("M130").Value = 3025
$ws.Range("N130").Value = -18290

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 60335.332
$ws.Range("I11").Value = 20998
$ws.Range("J11").Value = 80004
$ws.Range("K11").Value = 20998
$ws.Range("L11").Value = 80004
$ws.Range("M11").Value = -20859
$ws.Range("N11").Value = -80282

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1416.3572
$ws.Range("I46").Value = 1335.75
$ws.Range("J46").Value = 1900
$ws.Range("K46").Value = 1335.75
$ws.Range("L46").Value = 1900
$ws.Range("M46").Value = -1147.75
$ws.Range("N46").Value = -2276

$ws.Range("H55").Value = 587.3684
$ws.Range("I55").Value = 809.9
$ws.Range("J55").Value = 340.1111
$ws.Range("K55").Value = 809.9
$ws.Range("L55").Value = 340.1111
$ws.Range("M55").Value = -636.9
$ws.Range("N55").Value = -686.1111000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1688.2941
$ws.Range("I81").Value = 1515.4615
$ws.Range("J81").Value = 2250
$ws.Range("K81").Value = 3030.923
$ws.Range("L81").Value = 4500
$ws.Range("M81").Value = -1969.923
$ws.Range("N81").Value = -6622

$ws.Range("H84").Value = 1688.2941
$ws.Range("I84").Value = 1515.4615
$ws.Range("J84").Value = 2250
$ws.Range("K84").Value = 15154.615
$ws.Range("L84").Value = 22500
$ws.Range("M84").Value = -9850.614999999998
$ws.Range("N84").Value = -33108

